$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 1614.3889
$ws.Range("I19").Value = 586.7778
$ws.Range("J19").Value = 2642
$ws.Range("K19").Value = 586.7778
$ws.Range("L19").Value = 2642
$ws.Range("M19").Value = -411.7778
$ws.Range("N19").Value = -2992
$ws.Range("H80").Value = 3541.9375
$ws.Range("J80").Value = 3275.818
$ws.Range("L80").Value = 9827.454000000002
$ws.Range("N80").Value = -11823.454
$ws.Range("H83").Value = 3541.9375
$ws.Range("J83").Value = 3275.818
$ws.Range("L83").Value = 29482.362
$ws.Range("N83").Value = -39466.362
$ws.Range("H94").Value = 3000
$ws.Range("I94").Value = 3000
$ws.Range("K94").Value = 3000
$ws.Range("M94").Value = -2549
$ws.Range("H132").Value = 908.43634
$ws.Range("I132").Value = 770.0833
$ws.Range("J132").Value = 1857.1428
$ws.Range("K132").Value = 2310.2499
$ws.Range("L132").Value = 5571.428400000001
$ws.Range("M132").Value = 219.7501000000002
$ws.Range("N132").Value = -10631.4284
$ws.Range("H137").Value = 92956.27
$ws.Range("I137").Value = 3002
$ws.Range("J137").Value = 101951.7
$ws.Range("K137").Value = 9006
$ws.Range("L137").Value = 305855.1
$ws.Range("M137").Value = -6456
$ws.Range("N137").Value = -310955.1
$ws.Range("H138").Value = 1512.56
$ws.Range("I138").Value = 1239.125
$ws.Range("J138").Value = 1641.2354
$ws.Range("K138").Value = 3717.375
$ws.Range("L138").Value = 4923.706200000001
$ws.Range("M138").Value = 1422.625
$ws.Range("N138").Value = -15203.7062

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 293383.16
$ws.Range("I2").Value = 505789.38
$ws.Range("J2").Value = 1324.625
$ws.Range("K2").Value = 505789.38
$ws.Range("L2").Value = 1324.625
$ws.Range("M2").Value = -505676.38
$ws.Range("N2").Value = -1550.625
$ws.Range("H32").Value = 6184.28
$ws.Range("I32").Value = 3615.6782
$ws.Range("J32").Value = 23374.154
$ws.Range("K32").Value = 3615.6782
$ws.Range("L32").Value = 23374.154
$ws.Range("M32").Value = -3328.6782
$ws.Range("N32").Value = -23948.154
$ws.Range("H61").Value = 35540.25
$ws.Range("I61").Value = 45247.61
$ws.Range("K61").Value = 45247.61
$ws.Range("M61").Value = -45035.61
$ws.Range("H74").Value = 766.7742
$ws.Range("I74").Value = 766.7742
$ws.Range("K74").Value = 766.7742
$ws.Range("M74").Value = 107.2258
$ws.Range("H77").Value = 766.7742
$ws.Range("I77").Value = 766.7742
$ws.Range("K77").Value = 3833.871
$ws.Range("M77").Value = 534.1290000000004
$ws.Range("H116").Value = 293383.16
$ws.Range("I116").Value = 505789.38
$ws.Range("J116").Value = 1324.625
$ws.Range("K116").Value = 505789.38
$ws.Range("L116").Value = 1324.625
$ws.Range("M116").Value = -503495.38
$ws.Range("N116").Value = -5912.625
$ws.Range("H126").Value = 2997.5
$ws.Range("I126").Value = 2997.5
$ws.Range("K126").Value = 8992.5
$ws.Range("M126").Value = -6522.5
$ws.Range("H136").Value = 35540.25
$ws.Range("I136").Value = 45247.61
$ws.Range("K136").Value = 135742.83
$ws.Range("M136").Value = -133192.83

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 293383.16
$ws.Range("I3").Value = 505789.38
$ws.Range("J3").Value = 1324.625
$ws.Range("K3").Value = 505789.38
$ws.Range("L3").Value = 1324.625
$ws.Range("M3").Value = -505675.38
$ws.Range("N3").Value = -1552.625
$ws.Range("H20").Value = 4009.889
$ws.Range("I20").Value = 2899
$ws.Range("J20").Value = 4898.6
$ws.Range("K20").Value = 2899
$ws.Range("L20").Value = 4898.6
$ws.Range("M20").Value = -2652
$ws.Range("N20").Value = -5392.6
$ws.Range("H99").Value = 1550
$ws.Range("J99").Value = 2004
$ws.Range("L99").Value = 2004
$ws.Range("N99").Value = -5000
$ws.Range("H134").Value = 6018.2925
$ws.Range("I134").Value = 6368.3125
$ws.Range("J134").Value = 4773.778
$ws.Range("K134").Value = 19104.9375
$ws.Range("L134").Value = 14321.334
$ws.Range("M134").Value = -16569.9375
$ws.Range("N134").Value = -19391.334

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2667.111
$ws.Range("I31").Value = 1813.6666
$ws.Range("K31").Value = 1813.6666
$ws.Range("M31").Value = -1518.6666
$ws.Range("H34").Value = 2667.111
$ws.Range("I34").Value = 1813.6666
$ws.Range("K34").Value = 1813.6666
$ws.Range("M34").Value = -1611.6666
$ws.Range("H132").Value = 1415
$ws.Range("I132").Value = 864.1667
$ws.Range("J132").Value = 4059
$ws.Range("K132").Value = 2592.5001
$ws.Range("L132").Value = 12177
$ws.Range("M132").Value = -62.5001000000002
$ws.Range("N132").Value = -17237
$ws.Range("H134").Value = 1241.8096
$ws.Range("I134").Value = 1116.1471
$ws.Range("K134").Value = 3348.4413
$ws.Range("M134").Value = -813.4412999999995

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 125778.4
$ws.Range("I4").Value = 63851.875
$ws.Range("K4").Value = 191555.625
$ws.Range("M4").Value = -191443.625
$ws.Range("H33").Value = 166.76923
$ws.Range("I33").Value = 49.833332
$ws.Range("K33").Value = 298.999992
$ws.Range("M33").Value = -15.99999200000002
$ws.Range("H37").Value = 49999
$ws.Range("J37").Value = 49999
$ws.Range("L37").Value = 149997
$ws.Range("N37").Value = -150221
$ws.Range("H80").Value = 1694.4445
$ws.Range("H83").Value = 1694.4445
$ws.Range("H131").Value = 14719.759
$ws.Range("I131").Value = 398.8
$ws.Range("J131").Value = 16070.792
$ws.Range("K131").Value = 1196.4
$ws.Range("L131").Value = 48212.376
$ws.Range("M131").Value = 3843.6
$ws.Range("N131").Value = -58292.376

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H7").Value = 3870555.5
$ws.Range("J7").Value = 861666.3
$ws.Range("L7").Value = 861666.3
$ws.Range("N7").Value = -861890.3
$ws.Range("H8").Value = 3870555.5
$ws.Range("J8").Value = 861666.3
$ws.Range("L8").Value = 861666.3
$ws.Range("N8").Value = -861944.3
$ws.Range("H113").Value = 1574.75
$ws.Range("J113").Value = 1574.75
$ws.Range("L113").Value = 1574.75
$ws.Range("N113").Value = -5914.75
$ws.Range("H132").Value = 702907.8
$ws.Range("I132").Value = 839651.5600000001
$ws.Range("K132").Value = 2518954.68
$ws.Range("M132").Value = -2516424.68

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 2498.4211
$ws.Range("I7").Value = 1973.4546
$ws.Range("K7").Value = 1973.4546
$ws.Range("M7").Value = -1861.4546
$ws.Range("H61").Value = 2469.5557
$ws.Range("I61").Value = 2153.25
$ws.Range("K61").Value = 2153.25
$ws.Range("M61").Value = -1951.25
$ws.Range("H113").Value = 2469.5557
$ws.Range("I113").Value = 2153.25
$ws.Range("K113").Value = 2153.25
$ws.Range("M113").Value = 16.75
$ws.Range("H126").Value = 2498.4211
$ws.Range("I126").Value = 1973.4546
$ws.Range("K126").Value = 5920.3638
$ws.Range("M126").Value = -3450.3638
$ws.Range("H127").Value = 32972
$ws.Range("J127").Value = 32972
$ws.Range("L127").Value = 32972
$ws.Range("N127").Value = -42892
$ws.Range("H132").Value = 3564.7
$ws.Range("I132").Value = 1932.5
$ws.Range("K132").Value = 5797.5
$ws.Range("M132").Value = -3267.5
$ws.Range("H136").Value = 2791.2942
$ws.Range("I136").Value = 2776.8
$ws.Range("K136").Value = 8330.400000000001
$ws.Range("M136").Value = -5780.400000000001

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 566.1579
$ws.Range("I107").Value = 506.82352
$ws.Range("J107").Value = 1070.5
$ws.Range("K107").Value = 1520.47056
$ws.Range("L107").Value = 3211.5
$ws.Range("M107").Value = 399.52944
$ws.Range("N107").Value = -7051.5
$ws.Range("H122").Value = 64550.77
$ws.Range("I122").Value = 91434.44500000001
$ws.Range("J122").Value = 4062.5
$ws.Range("K122").Value = 274303.335
$ws.Range("L122").Value = 12187.5
$ws.Range("M122").Value = -271853.335
$ws.Range("N122").Value = -17087.5
$ws.Range("H132").Value = 1052.5063
$ws.Range("I132").Value = 962.9091
$ws.Range("J132").Value = 1507.3846
$ws.Range("K132").Value = 2888.7273
$ws.Range("L132").Value = 4522.1538
$ws.Range("M132").Value = -358.7273
$ws.Range("N132").Value = -9582.1538
